# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit refresh to Leve profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row40
$ws.Cells.Item(40, 8).Value = 2293.6365  # H40: was 2354.1667
$ws.Cells.Item(40, 9).Value = 2403.75  # I40: was 2164.2856
$ws.Cells.Item(40, 10).Value = 2000  # J40: was 2620
$ws.Cells.Item(40, 11).Value = 2403.75  # K40: was 2164.2856
$ws.Cells.Item(40, 12).Value = 2000  # L40: was 2620
$ws.Cells.Item(40, 13).Value = -2228.75  # M40: was -1989.2856
$ws.Cells.Item(40, 14).Value = -2350  # N40: was -2970

# ALC!row63
$ws.Cells.Item(63, 8).Value = 37950  # H63: was 38000
$ws.Cells.Item(63, 9).Value = 0  # I63: was 0
$ws.Cells.Item(63, 10).Value = 37950  # J63: was 38000
$ws.Cells.Item(63, 11).Value = 0  # K63: was 0
$ws.Cells.Item(63, 12).Value = 37950  # L63: was 38000
$ws.Cells.Item(63, 14).Value = -39198  # N63: was -39248

# ALC!row66
$ws.Cells.Item(66, 8).Value = 37950  # H66: was 38000
$ws.Cells.Item(66, 9).Value = 0  # I66: was 0
$ws.Cells.Item(66, 10).Value = 37950  # J66: was 38000
$ws.Cells.Item(66, 11).Value = 0  # K66: was 0
$ws.Cells.Item(66, 12).Value = 113850  # L66: was 114000
$ws.Cells.Item(66, 14).Value = -120090  # N66: was -120240

# ALC!row133
$ws.Cells.Item(133, 8).Value = 45000  # H133: was 39750
$ws.Cells.Item(133, 9).Value = 0  # I133: was 0
$ws.Cells.Item(133, 10).Value = 45000  # J133: was 39750
$ws.Cells.Item(133, 11).Value = 0  # K133: was 0
$ws.Cells.Item(133, 12).Value = 45000  # L133: was 39750
$ws.Cells.Item(133, 14).Value = -55120  # N133: was -49870

$ws = $wb.Worksheets.Item("ARM")
# ARM!row4
$ws.Cells.Item(4, 8).Value = 20220.4  # H4: was 100001
$ws.Cells.Item(4, 9).Value = 50051  # I4: was 100001
$ws.Cells.Item(4, 10).Value = 333.33334  # J4: was 0
$ws.Cells.Item(4, 11).Value = 50051  # K4: was 100001
$ws.Cells.Item(4, 12).Value = 333.33334  # L4: was 0
$ws.Cells.Item(4, 13).Value = -49935  # M4: was -99885
$ws.Cells.Item(4, 14).Value = -565.33334  # N4: was None

# ARM!row5
$ws.Cells.Item(5, 8).Value = 50082.332  # H5: was 60080.8
$ws.Cells.Item(5, 9).Value = 50073  # I5: was 66734
$ws.Cells.Item(5, 10).Value = 50101  # J5: was 50101
$ws.Cells.Item(5, 11).Value = 50073  # K5: was 66734
$ws.Cells.Item(5, 12).Value = 50101  # L5: was 50101
$ws.Cells.Item(5, 13).Value = -49961  # M5: was -66622

# ARM!row23
$ws.Cells.Item(23, 8).Value = 80006.5  # H23: was 35202.6
$ws.Cells.Item(23, 9).Value = 70006  # I23: was 70006
$ws.Cells.Item(23, 10).Value = 90007  # J23: was 26501.75
$ws.Cells.Item(23, 11).Value = 70006  # K23: was 70006
$ws.Cells.Item(23, 12).Value = 90007  # L23: was 26501.75
$ws.Cells.Item(23, 13).Value = -69747  # M23: was -69747
$ws.Cells.Item(23, 14).Value = -90525  # N23: was -27019.75

# ARM!row37
$ws.Cells.Item(37, 8).Value = 10000  # H37: was 11580
$ws.Cells.Item(37, 9).Value = 10000  # I37: was 11580
$ws.Cells.Item(37, 10).Value = 0  # J37: was 0
$ws.Cells.Item(37, 11).Value = 10000  # K37: was 11580
$ws.Cells.Item(37, 12).Value = 0  # L37: was 0
$ws.Cells.Item(37, 13).Value = -9727  # M37: was -11307

# ARM!row45
$ws.Cells.Item(45, 8).Value = 1683.7646  # H45: was 1179.3684
$ws.Cells.Item(45, 9).Value = 1147  # I45: was 1080.875
$ws.Cells.Item(45, 10).Value = 3428.25  # J45: was 1704.6666
$ws.Cells.Item(45, 11).Value = 1147  # K45: was 1080.875
$ws.Cells.Item(45, 12).Value = 3428.25  # L45: was 1704.6666
$ws.Cells.Item(45, 13).Value = -770  # M45: was -703.875
$ws.Cells.Item(45, 14).Value = -4182.25  # N45: was -2458.6666

# ARM!row55
$ws.Cells.Item(55, 8).Value = 0  # H55: was 44000
$ws.Cells.Item(55, 9).Value = 0  # I55: was 0
$ws.Cells.Item(55, 10).Value = 0  # J55: was 44000
$ws.Cells.Item(55, 11).Value = 0  # K55: was 0
$ws.Cells.Item(55, 12).Value = 0  # L55: was 44000
$ws.Cells.Item(55, 14).ClearContents()  # N55: was -44630, now absent

# ARM!row63
$ws.Cells.Item(63, 8).Value = 9237.714  # H63: was 13302.777
$ws.Cells.Item(63, 9).Value = 9717.615  # I63: was 13302.777
$ws.Cells.Item(63, 10).Value = 2999  # J63: was 0
$ws.Cells.Item(63, 11).Value = 9717.615  # K63: was 13302.777
$ws.Cells.Item(63, 12).Value = 2999  # L63: was 0
$ws.Cells.Item(63, 13).Value = -9031.615  # M63: was -12616.777
$ws.Cells.Item(63, 14).Value = -4371  # N63: was None

# ARM!row66
$ws.Cells.Item(66, 8).Value = 9237.714  # H66: was 13302.777
$ws.Cells.Item(66, 9).Value = 9717.615  # I66: was 13302.777
$ws.Cells.Item(66, 10).Value = 2999  # J66: was 0
$ws.Cells.Item(66, 11).Value = 48588.075  # K66: was 66513.88499999999
$ws.Cells.Item(66, 12).Value = 14995  # L66: was 0
$ws.Cells.Item(66, 13).Value = -45156.075  # M66: was -63081.88499999999
$ws.Cells.Item(66, 14).Value = -21859  # N66: was None

# ARM!row102
$ws.Cells.Item(102, 8).Value = 1500  # H102: was 1957.25
$ws.Cells.Item(102, 9).Value = 1500  # I102: was 1943
$ws.Cells.Item(102, 10).Value = 0  # J102: was 2000
$ws.Cells.Item(102, 11).Value = 1500  # K102: was 1943
$ws.Cells.Item(102, 12).Value = 0  # L102: was 2000
$ws.Cells.Item(102, 13).Value = 122  # M102: was -321
$ws.Cells.Item(102, 14).ClearContents()  # N102: was -5244, now absent

# ARM!row122
$ws.Cells.Item(122, 8).Value = 3080.2  # H122: was 3526.4707
$ws.Cells.Item(122, 9).Value = 2358.2942  # I122: was 2803.6
$ws.Cells.Item(122, 10).Value = 4614.25  # J122: was 4559.143
$ws.Cells.Item(122, 11).Value = 7074.882599999999  # K122: was 8410.799999999999
$ws.Cells.Item(122, 12).Value = 13842.75  # L122: was 13677.429
$ws.Cells.Item(122, 13).Value = -4624.882599999999  # M122: was -5960.799999999999
$ws.Cells.Item(122, 14).Value = -18742.75  # N122: was -18577.429

$ws = $wb.Worksheets.Item("BSM")
# BSM!row4
$ws.Cells.Item(4, 8).Value = 50082.332  # H4: was 60080.8
$ws.Cells.Item(4, 9).Value = 50073  # I4: was 66734
$ws.Cells.Item(4, 10).Value = 50101  # J4: was 50101
$ws.Cells.Item(4, 11).Value = 50073  # K4: was 66734
$ws.Cells.Item(4, 12).Value = 50101  # L4: was 50101
$ws.Cells.Item(4, 13).Value = -49958  # M4: was -66619

# BSM!row22
$ws.Cells.Item(22, 8).Value = 7796.923  # H22: was 7245.7144
$ws.Cells.Item(22, 9).Value = 8440.833000000001  # I22: was 8441.666999999999
$ws.Cells.Item(22, 10).Value = 70  # J22: was 70
$ws.Cells.Item(22, 11).Value = 8440.833000000001  # K22: was 8441.666999999999
$ws.Cells.Item(22, 12).Value = 70  # L22: was 70
$ws.Cells.Item(22, 13).Value = -8267.833000000001  # M22: was -8268.666999999999

# BSM!row94
$ws.Cells.Item(94, 8).Value = 704.8946999999999  # H94: was 800.1539
$ws.Cells.Item(94, 9).Value = 586  # I94: was 736.4286
$ws.Cells.Item(94, 10).Value = 868.375  # J94: was 874.5
$ws.Cells.Item(94, 11).Value = 586  # K94: was 736.4286
$ws.Cells.Item(94, 12).Value = 868.375  # L94: was 874.5
$ws.Cells.Item(94, 13).Value = -135  # M94: was -285.4286
$ws.Cells.Item(94, 14).Value = -1770.375  # N94: was -1776.5

# BSM!row107
$ws.Cells.Item(107, 8).Value = 1200  # H107: was 999.8
$ws.Cells.Item(107, 9).Value = 0  # I107: was 924.75
$ws.Cells.Item(107, 10).Value = 1200  # J107: was 1300
$ws.Cells.Item(107, 11).Value = 0  # K107: was 924.75
$ws.Cells.Item(107, 12).Value = 1200  # L107: was 1300
$ws.Cells.Item(107, 13).ClearContents()  # M107: was 995.25, now absent
$ws.Cells.Item(107, 14).Value = -5040  # N107: was -5140

$ws = $wb.Worksheets.Item("CRP")
# CRP!row22
$ws.Cells.Item(22, 8).Value = 1465.5555  # H22: was 2851
$ws.Cells.Item(22, 9).Value = 1611.25  # I22: was 5250
$ws.Cells.Item(22, 10).Value = 300  # J22: was 452
$ws.Cells.Item(22, 11).Value = 1611.25  # K22: was 5250
$ws.Cells.Item(22, 12).Value = 300  # L22: was 452
$ws.Cells.Item(22, 13).Value = -1261.25  # M22: was -4900
$ws.Cells.Item(22, 14).Value = -1000  # N22: was -1152

$ws = $wb.Worksheets.Item("CUL")
# CUL!row117
$ws.Cells.Item(117, 8).Value = 0  # H117: was 1589.2307
$ws.Cells.Item(117, 9).Value = 0  # I117: was 529
$ws.Cells.Item(117, 10).Value = 0  # J117: was 1677.5834
$ws.Cells.Item(117, 11).Value = 0  # K117: was 1587
$ws.Cells.Item(117, 12).Value = 0  # L117: was 5032.7502
$ws.Cells.Item(117, 13).ClearContents()  # M117: was 1855, now absent
$ws.Cells.Item(117, 14).ClearContents()  # N117: was -11916.7502, now absent

# CUL!row136
$ws.Cells.Item(136, 8).Value = 2312  # H136: was 2157.6667
$ws.Cells.Item(136, 9).Value = 2520  # I136: was 1989.875
$ws.Cells.Item(136, 10).Value = 2000  # J136: was 3500
$ws.Cells.Item(136, 11).Value = 7560  # K136: was 5969.625
$ws.Cells.Item(136, 12).Value = 6000  # L136: was 10500
$ws.Cells.Item(136, 13).Value = -2460  # M136: was -869.625
$ws.Cells.Item(136, 14).Value = -16200  # N136: was -20700

$ws = $wb.Worksheets.Item("GSM")
# GSM!row2
$ws.Cells.Item(2, 8).Value = 500  # H2: was 36.666668
$ws.Cells.Item(2, 9).Value = 0  # I2: was 10
$ws.Cells.Item(2, 10).Value = 500  # J2: was 50
$ws.Cells.Item(2, 11).Value = 0  # K2: was 10
$ws.Cells.Item(2, 12).Value = 500  # L2: was 50
$ws.Cells.Item(2, 13).ClearContents()  # M2: was 103, now absent
$ws.Cells.Item(2, 14).Value = -726  # N2: was -276

# GSM!row80
$ws.Cells.Item(80, 8).Value = 3258.25  # H80: was 4751.4
$ws.Cells.Item(80, 9).Value = 3166.5  # I80: was 6266.4
$ws.Cells.Item(80, 10).Value = 3350  # J80: was 3236.4
$ws.Cells.Item(80, 11).Value = 3166.5  # K80: was 6266.4
$ws.Cells.Item(80, 12).Value = 3350  # L80: was 3236.4
$ws.Cells.Item(80, 13).Value = -2168.5  # M80: was -5268.4
$ws.Cells.Item(80, 14).Value = -5346  # N80: was -5232.4

# GSM!row83
$ws.Cells.Item(83, 8).Value = 3258.25  # H83: was 4751.4
$ws.Cells.Item(83, 9).Value = 3166.5  # I83: was 6266.4
$ws.Cells.Item(83, 10).Value = 3350  # J83: was 3236.4
$ws.Cells.Item(83, 11).Value = 15832.5  # K83: was 31332
$ws.Cells.Item(83, 12).Value = 16750  # L83: was 16182
$ws.Cells.Item(83, 13).Value = -10840.5  # M83: was -26340
$ws.Cells.Item(83, 14).Value = -26734  # N83: was -26166

# GSM!row102
$ws.Cells.Item(102, 8).Value = 1750  # H102: was 1336.2354
$ws.Cells.Item(102, 9).Value = 1566.6666  # I102: was 1147.3846
$ws.Cells.Item(102, 10).Value = 1933.3334  # J102: was 1950
$ws.Cells.Item(102, 11).Value = 1566.6666  # K102: was 1147.3846
$ws.Cells.Item(102, 12).Value = 1933.3334  # L102: was 1950
$ws.Cells.Item(102, 13).Value = 55.33339999999998  # M102: was 474.6153999999999
$ws.Cells.Item(102, 14).Value = -5177.3334  # N102: was -5194

# GSM!row126
$ws.Cells.Item(126, 8).Value = 2048.2964  # H126: was 1765.975
$ws.Cells.Item(126, 9).Value = 1656.2778  # I126: was 1456.8148
$ws.Cells.Item(126, 10).Value = 2832.3333  # J126: was 2408.077
$ws.Cells.Item(126, 11).Value = 4968.8334  # K126: was 4370.4444
$ws.Cells.Item(126, 12).Value = 8496.999899999999  # L126: was 7224.231000000001
$ws.Cells.Item(126, 13).Value = -2498.8334  # M126: was -1900.4444
$ws.Cells.Item(126, 14).Value = -13436.9999  # N126: was -12164.231

$ws = $wb.Worksheets.Item("LTW")
# LTW!row22
$ws.Cells.Item(22, 8).Value = 1231.8077  # H22: was 1369.381
$ws.Cells.Item(22, 9).Value = 1001.2105  # I22: was 1097.0625
$ws.Cells.Item(22, 10).Value = 1857.7142  # J22: was 2240.8
$ws.Cells.Item(22, 11).Value = 1001.2105  # K22: was 1097.0625
$ws.Cells.Item(22, 12).Value = 1857.7142  # L22: was 2240.8
$ws.Cells.Item(22, 13).Value = -706.2105  # M22: was -802.0625
$ws.Cells.Item(22, 14).Value = -2447.7142  # N22: was -2830.8

# LTW!row27
$ws.Cells.Item(27, 8).Value = 1231.8077  # H27: was 1369.381
$ws.Cells.Item(27, 9).Value = 1001.2105  # I27: was 1097.0625
$ws.Cells.Item(27, 10).Value = 1857.7142  # J27: was 2240.8
$ws.Cells.Item(27, 11).Value = 1001.2105  # K27: was 1097.0625
$ws.Cells.Item(27, 12).Value = 1857.7142  # L27: was 2240.8
$ws.Cells.Item(27, 13).Value = -894.2105  # M27: was -990.0625
$ws.Cells.Item(27, 14).Value = -2071.7142  # N27: was -2454.8

# LTW!row46
$ws.Cells.Item(46, 8).Value = 63802.688  # H46: was 51218.3
$ws.Cells.Item(46, 9).Value = 143796.42  # I46: was 143714
$ws.Cells.Item(46, 10).Value = 1585.3334  # J46: was 1412.9231
$ws.Cells.Item(46, 11).Value = 143796.42  # K46: was 143714
$ws.Cells.Item(46, 12).Value = 1585.3334  # L46: was 1412.9231
$ws.Cells.Item(46, 13).Value = -143608.42  # M46: was -143526
$ws.Cells.Item(46, 14).Value = -1961.3334  # N46: was -1788.9231

# LTW!row55
$ws.Cells.Item(55, 8).Value = 224.21053  # H55: was 263.5625
$ws.Cells.Item(55, 9).Value = 239  # I55: was 280.15384
$ws.Cells.Item(55, 10).Value = 168.75  # J55: was 191.66667
$ws.Cells.Item(55, 11).Value = 239  # K55: was 280.15384
$ws.Cells.Item(55, 12).Value = 168.75  # L55: was 191.66667
$ws.Cells.Item(55, 13).Value = -66  # M55: was -107.15384
$ws.Cells.Item(55, 14).Value = -514.75  # N55: was -537.6666700000001

# LTW!row100
$ws.Cells.Item(100, 8).Value = 3350  # H100: was 3750
$ws.Cells.Item(100, 9).Value = 3350  # I100: was 3750
$ws.Cells.Item(100, 10).Value = 0  # J100: was 0
$ws.Cells.Item(100, 11).Value = 3350  # K100: was 3750
$ws.Cells.Item(100, 12).Value = 0  # L100: was 0
$ws.Cells.Item(100, 13).Value = -2809  # M100: was -3209

# LTW!row136
$ws.Cells.Item(136, 8).Value = 4839.8423  # H136: was 5682.357
$ws.Cells.Item(136, 9).Value = 4034.7693  # I136: was 4754.8
$ws.Cells.Item(136, 10).Value = 6584.1665  # J136: was 8001.25
$ws.Cells.Item(136, 11).Value = 12104.3079  # K136: was 14264.4
$ws.Cells.Item(136, 12).Value = 19752.4995  # L136: was 24003.75
$ws.Cells.Item(136, 13).Value = -9554.3079  # M136: was -11714.4
$ws.Cells.Item(136, 14).Value = -24852.4995  # N136: was -29103.75

$ws = $wb.Worksheets.Item("WVR")
# WVR!row96
$ws.Cells.Item(96, 8).Value = 1264.579  # H96: was 1502.5
$ws.Cells.Item(96, 9).Value = 1132.0834  # I96: was 1253
$ws.Cells.Item(96, 10).Value = 1491.7142  # J96: was 1752
$ws.Cells.Item(96, 11).Value = 1132.0834  # K96: was 1253
$ws.Cells.Item(96, 12).Value = 1491.7142  # L96: was 1752
$ws.Cells.Item(96, 13).Value = 240.9166  # M96: was 120
$ws.Cells.Item(96, 14).Value = -4237.7142  # N96: was -4498
